$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage for cells whose new
# value looks like a plain number (e.g. "25.20") so Excel does not silently
# coerce it to a Number and drop significant trailing zeros. The original
# cell style is restored immediately afterward so formatting is unaffected.
function Set-TextValue($cellRange, $value) {
    $origStyle = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = '60.213.36'
$ws.Range("E2").Value = '  -2.83%  '

# Row 3
$ws.Range("D3").Value = '2.385.18'
$ws.Range("E3").Value = '  -4.57%  '

# Row 5
Set-TextValue $ws.Range("D5") '538.91'
$ws.Range("E5").Value = '  -2.11%  '

# Row 6
Set-TextValue $ws.Range("D6") '139.81'
$ws.Range("E6").Value = '  -4.88%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("E8").Value = '  -5.82%  '

# Row 9
$ws.Range("D9").Value = '2.380.20'
$ws.Range("E9").Value = '  -4.72%  '

# Row 10
$ws.Range("E10").Value = '  -2.73%  '

# Row 11
$ws.Range("E11").Value = '  +0.01%  '

# Row 12
Set-TextValue $ws.Range("D12") '5.32'
$ws.Range("E12").Value = '  -1.18%  '

# Row 13
$ws.Range("E13").Value = '  -4.72%  '

# Row 14
Set-TextValue $ws.Range("D14") '25.20'
$ws.Range("E14").Value = '  -3.57%  '

# Row 15
$ws.Range("D15").Value = '2.816.37'
$ws.Range("E15").Value = '  -4.50%  '

# Row 16
$ws.Range("E16").Value = '  +0.32%  '

# Row 17
$ws.Range("D17").Value = '60.101.11'
$ws.Range("E17").Value = '  -2.82%  '

# Row 18
$ws.Range("D18").Value = '2.390.50'
$ws.Range("E18").Value = '  -4.17%  '

# Row 19
Set-TextValue $ws.Range("D19") '10.57'
$ws.Range("E19").Value = '  -4.73%  '

# Row 20
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D20") '4.05'
$ws.Range("E20").Value = '  -3.25%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D21") '6.68'
$ws.Range("E21").Value = '  -4.75%  '

# Row 22
Set-TextValue $ws.Range("D22") '313.45'
$ws.Range("E22").Value = '  -2.48%  '

# Row 23
$ws.Range("E23").Value = '  +0.02%  '

# Row 24
Set-TextValue $ws.Range("D24") '1.79'
$ws.Range("E24").Value = '  +3.33%  '

# Row 25
$ws.Range("E25").Value = '  -2.02%  '

# Row 26
$ws.Range("E26").Value = '  +0.52%  '

# Row 27
$ws.Range("D27").Value = '2.504.02'
$ws.Range("E27").Value = '  -4.56%  '

# Row 28
Set-TextValue $ws.Range("D28") '7.64'
$ws.Range("E28").Value = '  +0.15%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0909'
$ws.Range("E29").Value = '  -9.77%  '

# Row 30
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D30") '1.41'
$ws.Range("E30").Value = '  -5.03%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D31") '7.98'
$ws.Range("E31").Value = '  -4.79%  '

# Row 32
Set-TextValue $ws.Range("D32") '505.80'
$ws.Range("E32").Value = '  -5.41%  '

# Row 33
$ws.Range("E33").Value = '  -4.60%  '

# Row 34
$ws.Range("E34").Value = '  -3.90%  '

# Row 35
$ws.Range("E35").Value = '  -0.91%  '

# Row 36
$ws.Range("E36").Value = '  +0.10%  '

# Row 37
Set-TextValue $ws.Range("D37") '4.61'
$ws.Range("E37").Value = '  -5.22%  '

# Row 38
$ws.Range("E38").Value = '  -8.10%  '

# Row 39
$ws.Range("E39").Value = '  -1.74%  '

# Row 40
Set-TextValue $ws.Range("D40") '17.95'
$ws.Range("E40").Value = '  -2.91%  '

# Row 41
$ws.Range("E41").Value = '  +0.15%  '

# Row 42
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D42") '1.69'
$ws.Range("E42").Value = '  +0.02%  '

# Row 43
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D43") '136.87'
$ws.Range("E43").Value = '  -4.95%  '

# Row 44
Set-TextValue $ws.Range("D44") '40.28'
$ws.Range("E44").Value = '  -0.20%  '

# Row 45
Set-TextValue $ws.Range("D45") '2.16'
$ws.Range("E45").Value = '  -6.18%  '

# Row 46
Set-TextValue $ws.Range("D46") '139.56'
$ws.Range("E46").Value = '  -6.13%  '

# Row 47
Set-TextValue $ws.Range("D47") '3.51'
$ws.Range("E47").Value = '  -1.83%  '

# Row 48
Set-TextValue $ws.Range("D48") '20.15'
$ws.Range("E48").Value = '  -2.99%  '

# Row 49
Set-TextValue $ws.Range("D49") '0.0514'
$ws.Range("E49").Value = '  -3.99%  '

# Row 50
Set-TextValue $ws.Range("D50") '0.574'
$ws.Range("E50").Value = '  -2.39%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.0921'
$ws.Range("E51").Value = '  -2.62%  '
